# "Se guardan las clases reservadas" — save the reserved classes.
#
# GIMENEZ ZAIRA reserves Musculacion (LUNES 10:00) and Spinning (MARTES 08:00);
# Aerolocal's SABADO Localizada-style class also gets reserved. The C column on
# every class sheet is a 0/1 "reservado" flag; flip it to 1 for the classes
# that were just booked, and write the new student's booking sheet.

$wb = $excel.ActiveWorkbook

# --- GIMENEZ ZAIRA: record her two new reservations ------------------------
$zaira = $wb.Worksheets.Item("GIMENEZ ZAIRA")
$zaira.Range("A1").Value = "Musculacion"
$zaira.Range("B1").Value = "LUNES 10:00"
$zaira.Range("A2").Value = "Spinning"
$zaira.Range("B2").Value = "MARTES 08:00"

# --- Mark the reserved flag (column C) on each class schedule --------------

# Aerolocal: SABADO row (row 3) becomes reserved, and the sheet's selection
# ends up parked at H25.
$aerolocal = $wb.Worksheets.Item("Aerolocal")
$aerolocal.Range("C3").Value = 1
$aerolocal.Range("H25").Select()

# Funcional: MIERCOLES row (row 3) reservation re-saved.
$funcional = $wb.Worksheets.Item("Funcional")
$funcional.Range("C3").Value = 1

# Localizada: SABADO row (row 5) reservation re-saved.
$localizada = $wb.Worksheets.Item("Localizada")
$localizada.Range("C5").Value = 1

# Musculacion: LUNES 10:00 row (row 1) — Zaira's reservation.
$musculacion = $wb.Worksheets.Item("Musculacion")
$musculacion.Range("C1").Value = 1

# Spinning: MARTES row (row 2, Zaira's reservation) and MIERCOLES row 2 (row 4).
$spinning = $wb.Worksheets.Item("Spinning")
$spinning.Range("C2").Value = 1
$spinning.Range("C4").Value = 1

# Zumba: VIERNES row (row 5) reservation re-saved.
$zumba = $wb.Worksheets.Item("Zumba")
$zumba.Range("C5").Value = 1

# --- Restore focus to the sheet/cell the workbook had selected originally --
$sanchez = $wb.Worksheets.Item("SANCHEZ ROBERTO")
$sanchez.Activate()
$sanchez.Range("O29").Select()
